$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style/formatting from existing header cell H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-33
$data = @(
    @(7, 8),
    @(4, 5),
    @(5, 7),
    @(8, 9),
    @(7, 8),
    @(7, 8),
    @(6, 9),
    @(7, 8),
    @(1, 3),
    @(9, 9),
    @(7, 8),
    @(6, 7),
    @(6, 7),
    @(8, 9),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(7, 8),
    @(6, 7),
    @(6, 7),
    @(6, 8),
    @(7, 8),
    @(3, 6),
    @(1, 3),
    @(6, 7),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(5, 6),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
